# Weekly refresh of the price series: a new weekly reading is inserted at
# row 23 (pushing every existing record for this market down by one row),
# so the report now spans rows 2-82 instead of 2-81.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("23:23").Insert()

$ws.Range("A23").Value2 = 3
$ws.Range("B23").Value2 = "Femacal de La Calera"
$ws.Range("C23").Value2 = "Coquimbo"
$ws.Range("D23").Value2 = 44804
$ws.Range("E23").Value2 = 5
$ws.Range("F23").Value2 = 100112035
$ws.Range("G23").Value2 = "Bruselas (repollito)"
$ws.Range("H23").Value2 = "Sin especificar"
$ws.Range("I23").Value2 = "Primera"
$ws.Range("J23").Value2 = 50
$ws.Range("K23").Value2 = 15000
$ws.Range("L23").Value2 = 15000
$ws.Range("M23").Value2 = 15000
$ws.Range("N23").Value2 = "$/malla 15 kilos"
$ws.Range("O23").Value2 = "Provincia de Quillota"
$ws.Range("P23").Value2 = 1000
$ws.Range("Q23").Value2 = 15
$ws.Range("R23").Value2 = "Hortaliza"
